# The notebook that generates nlp/gen_goods.xlsx was re-run on updated data.
# The "Counts" table in column A/B is a word -> frequency table, sorted by
# descending frequency (column B is unchanged by this re-run). Within several
# tied-frequency groups (9, 5, 4, 2, 1) the tie-break order produced by the
# new run differs from the previous run, so the *word* shown in a handful of
# rows changes while its paired count in column B stays exactly the same.
#
# Below we just overwrite column A for the rows whose word changed, leaving
# every other cell (including all of column B) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# frequency-9 group (rows 17-18): order flips
$ws.Range("A17").Value = "деревенский товар"
$ws.Range("A18").Value = "серебреный товар"

# frequency-5 group (rows 24-27): new tie-break order
$ws.Range("A24").Value = "пушной товар"
$ws.Range("A25").Value = "нужный товар"
$ws.Range("A27").Value = "щепетильный товар"

# frequency-4 group (rows 28-32): new tie-break order
$ws.Range("A28").Value = "питейный припасы"
$ws.Range("A29").Value = "суровский товар"
$ws.Range("A30").Value = "медный товар"
$ws.Range("A31").Value = "внутренний товар"
$ws.Range("A32").Value = "недорогой товар"

# frequency-2 group (rows 36 & 38): order flips
$ws.Range("A36").Value = "купецкий товар"
$ws.Range("A38").Value = "галантерейный товар"

# frequency-1 group (rows 39-42): new tie-break order
$ws.Range("A39").Value = "меховой товар"
$ws.Range("A40").Value = "домовый товар"
$ws.Range("A41").Value = "харчевой припасы"
$ws.Range("A42").Value = "надлежащий товар"
$ws.Range("A43").Value = "рукодельный товар"
